# Commit: "Update data loading and add gitignore"
#
# Semantic changes applied here:
#   1. Rename the worksheet "Sheet2" -> "ITSS_Data".
#   2. Clear out any (now unused) conditional-formatting rules / their
#      differential styles (dxfs) on the sheet, mirroring the source
#      workbook no longer referencing the old red/green dxf pairs that
#      used to back the aging-bucket conditional formatting.
#
# Note: the workbook's `xr:revisionPtr` documentId GUID is regenerated by
# Excel itself on every save and isn't something a script should set
# explicitly, so it is intentionally left alone here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet.
$ws.Name = "ITSS_Data"

# 2. Drop every conditional-formatting rule left on the sheet so the
#    differential formats (dxfs) backing them are no longer referenced.
$ws.Cells.FormatConditions.Delete()
